$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 46: "article 45 is live"
$ws.Range("A46").Value = 45

$ws.Range("B46").Value = 43875
$ws.Range("B46").NumberFormat = "d-mmm-yy"

$ws.Range("C46").Value = @'
Surah Al Nisa, 30 - 36
'@

$ws.Range("F46").Value = @'
Justifying anger, Hitting women, Islamic society, Marriage, Responsibility of men and women
'@

$ws.Range("D46").Value = @'
h1: Hitting a woman in Islam is allowed.. but
quote: Men are in charge of women by [right of] what Allah has given one over the other and what they spend [for maintenance] from their wealth. So <span class=lavendar>righteous women are devoutly obedient, guarding in [the husband's] absence what Allah would have them guard</span>. But those [wives] from whom you fear arrogance - <span class=lavendar>[first] advise them</span>; [then if they persist], <span class=lavendar>forsake them in bed</span>; and [finally], <span class=lavendar>strike them</span>. But if they obey you [once more], <span class=lavendar>seek no means against them</span>. Indeed, Allah is ever Exalted and Grand. <br> - Surah Al Nosa verse 34
p: This is the highlighted version of Quran. Where I have selectively picked up one verse to justify <b>hitting the women</b>. It is allowed but Quran takes to this part of the verse gradually. 
p.b-left: <b>Verse 1.</b> Connect strongly to the <span class=lavendar>bonds of relationships</span>.
p.b-left: <b>Verse 3.</b> Give<span class=lavendar> justice to women </span>who are slaves or orphans.
p.b-left: <b>Verse 4. </b><span class=lavendar>Protect their property</span> in orphanage until they grow old.
p.b-left: <b>Verse 11. </b><span class=lavendar>Give her </span>1/6th to 1/8th of your wealth after you die.
p.b-left: <b>Verse 15. </b>If she commits adultery, bring 4 x witnesses, confine her and wait until she repents. After that, <span class=lavendar>free her</span>.
p.b-left: <b>Verse 17. </b><span class=lavendar>Do not make life difficult</span> for her.
p.b-left: <b>Verse 20. </b>If you gave her anything in marriage, <span class=lavendar>do not get it back</span> in divorce.
p.b-left: <b>Verse 24. </b>Marry the war prisoners, to <span class=lavendar>make them equal contributors </span>towards a healthy society.
p.b-left: <b>Verse 28. </b><span class=lavendar>Marriage is made easier </span>for us to keep the women protected.
p.b-left: <b>Verse 34. </b><span class=note>as mentioned in above verse</span>, to keep the marriage intact, advise women, forsake them in bed and strike them. If they obey, <span class=lavendar>seek no harm against them</span>.
h3: Linking the verses
p: The interlinking of these verses are emphasising on a united family system. The concept of relationship between man and woman is both romantic and purpose-built. It is also the most emotional feeling. Movies made on family life; <b>The blind side</b>, <b>The pursuit of happiness</b> and <b>The money ball</b> exploit the relationship within a family to bring in the emotional part. It is the fundamental fabric of a society. <b>Absence of a perfectly balanced family system gives spark to an imbalanced life</b>. 
p: The moment a broken woman is brought into custody through marriage, the responsibility of providing them shelter is 100 %. But if they end up violating the laid boundaries, the punishment on them is reduced to 50 %. We get angry with them 50 % lesser, get annoyed from them 50 % less and we love them 50 % more.
h3: Marriage is better than routine affair
p: Routine affairs between men and women do not bind us to helping each other. We still can untie, take a break and get away living in with somebody else. The focus towards building a better community needs us to marry earlier, sooner and blunt. 13 years old can marry 13 years old girls. Children can live a life of normalcy even after they are married earlier. 
p: <b>Marriage can be taken as just another subject in their school. </b>The subject that teaches bold realities of life and brings responsibility and focus to its best students.
h3: Hitting the women
p: “Hitting wives” is not needed. It comes after we start building a better society. It is only there in a <b>fundamentalists life</b> who picks a part of Quran to justify his anger. Anger is wrong. Hitting is right after anger is controlled. And after the anger is handled properly, healing broken relations becomes easier. Husbands become sober and they attract wives towards them. They start acting mature and truly embrace the status of a leader. They act leader-ly to inspire their families. 
p: Family is a fun loving fabric of our society. It needs to be healed.<b> Directly reaching the ‘hitting’ part is only anger justified</b>. It is not Islam and it is not human.
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>. Subscribe below to receive updates.
'@

$ws.Range("E46").Value = "Qasim Ali"

$ws.Rows.Item(46).RowHeight = 409.6

$ws.Range("D46").Select()
